$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$ws.Range("D2").Value = '68.896.29'
$ws.Range("D3").Value = '3.736.76'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("E4").Value = '  -0.03%  '
Set-TextValue "D5" '600.57'
$ws.Range("E5").Value = '  +0.05%  '
Set-TextValue "D6" '165.24'
$ws.Range("E6").Value = '  -2.21%  '
$ws.Range("D7").Value = '3.732.62'
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("E10").Value = '  +4.09%  '
Set-TextValue "D11" '6.41'
$ws.Range("E11").Value = '  +1.24%  '
$ws.Range("E12").Value = '  -0.38%  '
Set-TextValue "D13" '37.72'
$ws.Range("E13").Value = '  -0.94%  '
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("D15").Value = '4.365.10'
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = '3.737.88'
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").Value = '68.968.94'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("E18").Value = '  +2.24%  '
Set-TextValue "D19" '17.61'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("E20").Value = '  -2.04%  '
Set-TextValue "D21" '11.12'
$ws.Range("E21").Value = '  +5.16%  '
Set-TextValue "D22" '491.35'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("E23").Value = '  +0.09%  '
Set-TextValue "D24" '84.51'
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("E25").Value = '  +3.17%  '
$ws.Range("E26").Value = '  -1.76%  '
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D30" '8.21'
$ws.Range("E30").Value = '  +3.53%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D31" '2.96'
$ws.Range("E31").Value = '  +0.04%  '
Set-TextValue "D32" '2.42'
$ws.Range("E32").Value = '  -5.85%  '
$ws.Range("D33").Value = '3.885.11'
$ws.Range("E33").Value = '  +0.45%  '
Set-TextValue "D34" '31.51'
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("D35").Value = '3.673.71'
$ws.Range("E35").Value = '  +0.47%  '
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("E37").Value = '  +1.62%  '
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("E39").Value = '  +5.00%  '
Set-TextValue "D40" '3.17'
$ws.Range("E40").Value = '  +10.24%  '
$ws.Range("E42").Value = '  -0.18%  '
Set-TextValue "D43" '48.55'
$ws.Range("E43").Value = '  -0.65%  '
$ws.Range("E44").Value = '  +0.19%  '
Set-TextValue "D45" '423.59'
$ws.Range("E45").Value = '  -3.09%  '
Set-TextValue "D46" '8.43'
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("E47").Value = '  +0.01%  '
Set-TextValue "D48" '39.86'
$ws.Range("E48").Value = '  -1.75%  '
Set-TextValue "D49" '141.03'
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("D50").Value = '2.775.75'
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("E51").Value = '  +5.69%  '
